# Samples.xlsx update — re-brand the demo inventory data:
#   - Header rename: CDU1/CDU2 -> C1/C2
#   - Grid column:   "1S" -> "Rack1"
#   - Hostnames:     "hpcs-*" -> "hpSample-*"
#   - Management IPs (10.74.189.x / 10.74.188.x) -> 192.167.189.x / 192.167.188.x
#   - Port/Lom IPs   (10.9.24.x)                 -> 192.168.24.x
#   - Widen column A, move selection to D13 and drop the saved scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row renames (E1 = CDU1 -> C1, G1 = CDU2 -> C2) ---
$ws.Range("E1").Value = "C1"
$ws.Range("G1").Value = "C2"

# --- Data rows: walk every row in the table and rewrite the values in place ---
for ($r = 3; $r -le 52; $r++) {

    # Column C (Grid): "1S" -> "Rack1"
    $grid = $ws.Cells.Item($r, 3).Value2
    if ($grid -eq "1S") {
        $ws.Cells.Item($r, 3).Value = "Rack1"
    }

    # Column H (Hostname): "hpcs-..." -> "hpSample-..."
    $hostname = $ws.Cells.Item($r, 8).Value2
    if ($hostname -ne $null -and $hostname -ne "" -and $hostname.StartsWith("hpcs-")) {
        $ws.Cells.Item($r, 8).Value = $hostname.Replace("hpcs-", "hpSample-")
    }

    # Column K (IPv4): re-home both the management and the LOM/port subnets
    $ip = $ws.Cells.Item($r, 11).Value2
    if ($ip -ne $null -and $ip -ne "") {
        $newIp = $ip.Replace("10.74.189", "192.167.189").Replace("10.74.188", "192.167.188").Replace("10.9.24", "192.168.24")
        if ($newIp -ne $ip) {
            $ws.Cells.Item($r, 11).Value = $newIp
        }
    }
}

# --- Column A width widened to fit the (now longer) formula results ---
$ws.Columns.Item(1).ColumnWidth = 38.88

# --- Selection moves to D13, and the saved top-left scroll anchor is cleared ---
$ws.Range("D13").Select()
